$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 2).Value = 0.009057032445356359
$ws.Cells.Item(2, 3).Value = 0.2371322208235792
$ws.Cells.Item(3, 2).Value = -0.0707973000599024
$ws.Cells.Item(3, 3).Value = 0.3409380406977133
$ws.Cells.Item(4, 2).Value = 0.0172660329800893
$ws.Cells.Item(4, 3).Value = 0.3563682826067976
$ws.Cells.Item(5, 2).Value = 0.1842103288777725
$ws.Cells.Item(5, 3).Value = 0.06221521385363245
$ws.Cells.Item(6, 2).Value = 0.6674371044632009
$ws.Cells.Item(6, 3).Value = 0.1295006050101185
$ws.Cells.Item(7, 2).Value = 0.3709289385122904
$ws.Cells.Item(7, 3).Value = 0.1580800546445858
$ws.Cells.Item(8, 2).Value = -0.08485488332187253
$ws.Cells.Item(8, 3).Value = 0.40149525222309
$ws.Cells.Item(9, 2).Value = -0.3803102403821008
$ws.Cells.Item(9, 3).Value = 0.310282407843002
$ws.Cells.Item(10, 2).Value = -0.1370945519480693
$ws.Cells.Item(10, 3).Value = 0.01737396179205468
$ws.Cells.Item(11, 2).Value = 0.2654318459850528
$ws.Cells.Item(11, 3).Value = 0.3833718108625055
$ws.Cells.Item(12, 2).Value = 0.1925992380236453
$ws.Cells.Item(12, 3).Value = 0.1410302399371223
$ws.Cells.Item(13, 2).Value = 0.2671110743073307
$ws.Cells.Item(13, 3).Value = -0.1864056215647717
$ws.Cells.Item(14, 2).Value = 0.1524756792599838
$ws.Cells.Item(14, 3).Value = -0.3977223039149543
$ws.Cells.Item(15, 2).Value = 0.01643676288178535
$ws.Cells.Item(15, 3).Value = 0.1562223533630081
$ws.Cells.Item(16, 2).Value = -0.06110911106146517
$ws.Cells.Item(16, 3).Value = 0.1128511053235028
$ws.Cells.Item(17, 2).Value = 0.0288488885040994
$ws.Cells.Item(17, 3).Value = -0.01381828810323762
